$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 38640
$ws.Range("P3").Value = 8242.2
$ws.Range("U3").Value = "8242.2 (±444.2)"
$ws.Range("X3").Value = 2019.3
$ws.Range("Y3").Value = 108.8
$ws.Range("Z3").Value = "2019.3(±108.8)"

# Row 4
$ws.Range("D4").Value = 25542
$ws.Range("P4").Value = 2641.4
$ws.Range("S4").Value = 2.1
$ws.Range("U4").Value = "2641.4 (±452.1)"
$ws.Range("V4").Value = "11.5% (±2.1%)"
$ws.Range("X4").Value = 1196.2
$ws.Range("Y4").Value = 204.7
$ws.Range("Z4").Value = "1196.2(±204.7)"

# Row 5
$ws.Range("D5").Value = 15450
$ws.Range("P5").Value = 2013
$ws.Range("U5").Value = "2013.0 (±135.2)"
$ws.Range("X5").Value = 1376.7
$ws.Range("Y5").Value = 92.4
$ws.Range("Z5").Value = "1376.7(±92.4)"

# Row 8
$ws.Range("D8").Value = 4636
$ws.Range("P8").Value = 452
$ws.Range("U8").Value = "452.0 (±87.0)"
$ws.Range("X8").Value = 812.9
$ws.Range("Y8").Value = 156.5
$ws.Range("Z8").Value = "812.9(±156.5)"

# Row 10
$ws.Range("D10").Value = 201720
$ws.Range("P10").Value = 26978.6
$ws.Range("U10").Value = "26978.6 (±1669.3)"
$ws.Range("X10").Value = 1039.5
$ws.Range("Y10").Value = 64.3
$ws.Range("Z10").Value = "1039.5(±64.3)"

# Row 15
$ws.Range("D15").Value = 11633
$ws.Range("P15").Value = 1511
$ws.Range("R15").Value = 14.9
$ws.Range("S15").Value = 2.2
$ws.Range("U15").Value = "1511.0 (±204.2)"
$ws.Range("V15").Value = "14.9% (±2.2%)"
$ws.Range("X15").Value = 1280.7
$ws.Range("Y15").Value = 173
$ws.Range("Z15").Value = "1280.7(±173.0)"

# Row 18
$ws.Range("D18").Value = 1457
$ws.Range("P18").Value = 215
$ws.Range("R18").Value = 17.3
$ws.Range("U18").Value = "215.0 (±46.5)"
$ws.Range("V18").Value = "17.3% (±4.2%)"
$ws.Range("X18").Value = 1767.2
$ws.Range("Z18").Value = "1767.2(±382.2)"

# Row 22
$ws.Range("D22").Value = 36588
$ws.Range("P22").Value = 6336.6
$ws.Range("R22").Value = 20.9
$ws.Range("U22").Value = "6336.6 (±707.2)"
$ws.Range("V22").Value = "20.9% (±2.7%)"
$ws.Range("X22").Value = 1492.5
$ws.Range("Y22").Value = 166.5
$ws.Range("Z22").Value = "1492.5(±166.5)"

# Row 24
$ws.Range("D24").Value = 12830
$ws.Range("P24").Value = 1497.2
$ws.Range("R24").Value = 13.2
$ws.Range("U24").Value = "1497.2 (±199.3)"
$ws.Range("V24").Value = "13.2% (±1.9%)"
$ws.Range("X24").Value = 1199.7
$ws.Range("Z24").Value = "1199.7(±159.7)"

# Row 29
$ws.Range("D29").Value = 26749
$ws.Range("P29").Value = 6199.4
$ws.Range("U29").Value = "6199.4 (±289.3)"
$ws.Range("X29").Value = 2596.2
$ws.Range("Y29").Value = 121.1
$ws.Range("Z29").Value = "2596.2(±121.1)"

# Row 31
$ws.Range("D31").Value = 9093
$ws.Range("P31").Value = 1608.2
$ws.Range("R31").Value = 21.5
$ws.Range("U31").Value = "1608.2 (±111.7)"
$ws.Range("V31").Value = "21.5% (±1.8%)"
$ws.Range("X31").Value = 2252
$ws.Range("Z31").Value = "2252.0(±156.4)"

# Row 36
$ws.Range("D36").Value = 140519
$ws.Range("P36").Value = 21968.4
$ws.Range("U36").Value = "21968.4 (±1234.6)"
$ws.Range("X36").Value = 1528.7
$ws.Range("Z36").Value = "1528.7(±85.9)"

# Row 41
$ws.Range("D41").Value = 5482
$ws.Range("P41").Value = 879.4
$ws.Range("R41").Value = 19.1
$ws.Range("U41").Value = "879.4 (±40.2)"
$ws.Range("V41").Value = "19.1% (±1.0%)"
$ws.Range("X41").Value = 2029.7
$ws.Range("Y41").Value = 92.8
$ws.Range("Z41").Value = "2029.7(±92.8)"

# Row 42
$ws.Range("D42").Value = 942
$ws.Range("P42").Value = 188
$ws.Range("R42").Value = 24.9
$ws.Range("S42").Value = 4.3
$ws.Range("U42").Value = "188.0 (±27.3)"
$ws.Range("V42").Value = "24.9% (±4.3%)"
$ws.Range("X42").Value = 2087.5
$ws.Range("Z42").Value = "2087.5(±303.1)"

# Row 44
$ws.Range("D44").Value = 1098
$ws.Range("P44").Value = 239.4
$ws.Range("R44").Value = 27.9
$ws.Range("S44").Value = 5
$ws.Range("U44").Value = "239.4 (±35.0)"
$ws.Range("V44").Value = "27.9% (±5.0%)"
$ws.Range("X44").Value = 2965.4
$ws.Range("Y44").Value = 433.6
$ws.Range("Z44").Value = "2965.4(±433.6)"

# Row 48
$ws.Range("D48").Value = 25153
$ws.Range("P48").Value = 4220.6
$ws.Range("R48").Value = 20.2
$ws.Range("S48").Value = 3.5
$ws.Range("U48").Value = "4220.6 (±613.3)"
$ws.Range("V48").Value = "20.2% (±3.5%)"
$ws.Range("X48").Value = 1781.7
$ws.Range("Z48").Value = "1781.7(±258.9)"

# Row 52
$ws.Range("D52").Value = 111388
$ws.Range("P52").Value = 23919
$ws.Range("U52").Value = "23919.0 (±651.6)"
$ws.Range("X52").Value = 2238.5
$ws.Range("Y52").Value = 61
$ws.Range("Z52").Value = "2238.5(±61.0)"

# Row 55
$ws.Range("D55").Value = 65389
$ws.Range("P55").Value = 14441.6
$ws.Range("U55").Value = "14441.6 (±653.5)"
$ws.Range("X55").Value = 2232.2
$ws.Range("Z55").Value = "2232.2(±101.0)"

# Row 56
$ws.Range("D56").Value = 42830
$ws.Range("P56").Value = 5198.8
$ws.Range("U56").Value = "5198.8 (±757.3)"
$ws.Range("X56").Value = 1535.3
$ws.Range("Y56").Value = 223.7
$ws.Range("Z56").Value = "1535.3(±223.7)"

# Row 57
$ws.Range("D57").Value = 24543
$ws.Range("P57").Value = 3621.2
$ws.Range("U57").Value = "3621.2 (±214.7)"
$ws.Range("X57").Value = 1663.9
$ws.Range("Y57").Value = 98.7
$ws.Range("Z57").Value = "1663.9(±98.7)"

# Row 60
$ws.Range("D60").Value = 6692
$ws.Range("P60").Value = 720.8
$ws.Range("U60").Value = "720.8 (±148.8)"
$ws.Range("X60").Value = 964.7
$ws.Range("Y60").Value = 199.2
$ws.Range("Z60").Value = "964.7(±199.2)"

# Row 62
$ws.Range("D62").Value = 342239
$ws.Range("P62").Value = 48947
$ws.Range("S62").Value = 1.1
$ws.Range("U62").Value = "48947.0 (±2877.5)"
$ws.Range("V62").Value = "16.7% (±1.1%)"
$ws.Range("X62").Value = 1213.8
$ws.Range("Z62").Value = "1213.8(±71.4)"

# Row 67
$ws.Range("D67").Value = 17115
$ws.Range("P67").Value = 2390.4
$ws.Range("R67").Value = 16.2
$ws.Range("S67").Value = 1.8
$ws.Range("U67").Value = "2390.4 (±233.2)"
$ws.Range("V67").Value = "16.2% (±1.8%)"
$ws.Range("X67").Value = 1481.9
$ws.Range("Z67").Value = "1481.9(±144.5)"

# Row 68
$ws.Range("D68").Value = 2208
$ws.Range("P68").Value = 378.8
$ws.Range("R68").Value = 20.7
$ws.Range("U68").Value = "378.8 (±63.3)"
$ws.Range("V68").Value = "20.7% (±4.0%)"
$ws.Range("X68").Value = 1564.4
$ws.Range("Y68").Value = 261.5
$ws.Range("Z68").Value = "1564.4(±261.5)"

# Row 70
$ws.Range("D70").Value = 2555
$ws.Range("P70").Value = 454.4
$ws.Range("R70").Value = 21.6
$ws.Range("U70").Value = "454.4 (±72.2)"
$ws.Range("V70").Value = "21.6% (±4.0%)"
$ws.Range("X70").Value = 2245.2
$ws.Range("Z70").Value = "2245.2(±356.7)"

# Row 74
$ws.Range("D74").Value = 61741
$ws.Range("P74").Value = 10557.2
$ws.Range("R74").Value = 20.6
$ws.Range("S74").Value = 2.9
$ws.Range("U74").Value = "10557.2 (±1289.2)"
$ws.Range("V74").Value = "20.6% (±2.9%)"
$ws.Range("X74").Value = 1596.1
$ws.Range("Z74").Value = "1596.1(±194.9)"

# Row 76
$ws.Range("D76").Value = 19896
$ws.Range("P76").Value = 2732.2
$ws.Range("S76").Value = 1.9
$ws.Range("U76").Value = "2732.2 (±293.1)"
$ws.Range("V76").Value = "15.9% (±1.9%)"
$ws.Range("X76").Value = 1522
$ws.Range("Z76").Value = "1522.0(±163.3)"

# Row 78
$ws.Range("D78").Value = 263609
$ws.Range("P78").Value = 54899
$ws.Range("U78").Value = "54899.0 (±1249.5)"
$ws.Range("X78").Value = 1905.6
$ws.Range("Z78").Value = "1905.6(±43.4)"
